$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 29; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 45574
}
